{"js": "// The document ends with an empty trailing paragraph. Fill it with the\n// first new brainstorming bullet, then add a second bullet as a brand\n// new paragraph right after it.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\n\nlastParagraph.insertText(\n  \"Look at disparate impact (broad groups) and how it changes with each iteration (graph this)\",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n\nlastParagraph.insertParagraph(\n  \"Look at intersectionality and how it changes with each iteration\",\n  Word.InsertLocation.after\n);\nawait context.sync();\n", "ps1": "# Add two brainstorming bullets to the end of the document:\n#   1) fill the trailing empty paragraph with the first line\n#   2) insert a brand-new paragraph after it with the second line\n$d = $word.ActiveDocument\n\n$lastParagraph = $d.Paragraphs.Item($d.Paragraphs.Count)\n$lastParagraph.Range.Text = \"Look at disparate impact (broad groups) and how it changes with each iteration (graph this)\"\n\n$lastParagraph = $d.Paragraphs.Item($d.Paragraphs.Count)\n$lastParagraph.Range.InsertParagraphAfter()\n\n$newParagraph = $d.Paragraphs.Item($d.Paragraphs.Count)\n$newParagraph.Range.Text = \"Look at intersectionality and how it changes with each iteration\"\n"}
